# distribution_circuits.xlsx update
# "Getting a few more of the circuits up and running (Ckt 5 done; first part
#  of J1; 8500 halfway through taps). Some small bug fixes."

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 8 "8500 node": taps (Regs) progress, 3 -> 12; TODO severity bumped to HARD
$ws.Range("I8").Value = 12
$ws.Range("P8").Value = "HARD"

# --- Row 11 "Ckt5": HC Calcs finished (done) -> mark X, drop it from the TODO list
$ws.Range("N11").Value = "X"
$ws.Range("P11").ClearContents()
$ws.Range("Q11").ClearContents()

# --- Row 14 "J1": first part done -> Fixed model implemented
$ws.Range("L14").Value = "X"

# --- Small bug fix: stray marker value left in R6
$ws.Range("R6").Value = "  "

# Restore cursor/selection position as left by the author
$ws.Range("Q5").Select() | Out-Null

$wb.Save()
